# Apply updated cryptocurrency price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text (may look numeric); force text storage so
# formatting such as trailing zeros / multi-dot thousands separators is kept.
$priceCells = @(
    @{ Ref = 'D2'; Value = '65.761.21' },
    @{ Ref = 'D3'; Value = '2.678.97' },
    @{ Ref = 'D5'; Value = '600.47' },
    @{ Ref = 'D6'; Value = '156.97' },
    @{ Ref = 'D13'; Value = '29.37' },
    @{ Ref = 'D15'; Value = '3.159.42' },
    @{ Ref = 'D16'; Value = '66.413.61' },
    @{ Ref = 'D17'; Value = '2.675.69' },
    @{ Ref = 'D18'; Value = '12.93' },
    @{ Ref = 'D21'; Value = '352.71' },
    @{ Ref = 'D25'; Value = '9.65' },
    @{ Ref = 'D29'; Value = '8.04' },
    @{ Ref = 'D31'; Value = '529.89' },
    @{ Ref = 'D34'; Value = '6.48' },
    @{ Ref = 'D35'; Value = '5.50' },
    @{ Ref = 'D36'; Value = '20.76' },
    @{ Ref = 'D37'; Value = '0.424' },
    @{ Ref = 'D38'; Value = '1.00' },
    @{ Ref = 'D39'; Value = '158.06' },
    @{ Ref = 'D40'; Value = '1.95' },
    @{ Ref = 'D42'; Value = '164.39' },
    @{ Ref = 'D51'; Value = '20.23' }
)

foreach ($item in $priceCells) {
    $cell = $ws.Range($item.Ref)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Value
    $cell.Style = "Normal"
}

# Columns B, C and E hold plain text (coin name, link, volume %) and can
# be assigned directly.
$textCells = @(
    @{ Ref = 'E2'; Value = '  -0.02%  ' },
    @{ Ref = 'E4'; Value = '  +0.03%  ' },
    @{ Ref = 'E5'; Value = '  -0.84%  ' },
    @{ Ref = 'E6'; Value = '  -0.51%  ' },
    @{ Ref = 'E7'; Value = '  +0.01%  ' },
    @{ Ref = 'E8'; Value = '  +5.97%  ' },
    @{ Ref = 'E9'; Value = '  +5.67%  ' },
    @{ Ref = 'E10'; Value = '  -0.05%  ' },
    @{ Ref = 'E11'; Value = '  -2.54%  ' },
    @{ Ref = 'E12'; Value = '  +0.01%  ' },
    @{ Ref = 'E13'; Value = '  -2.43%  ' },
    @{ Ref = 'E14'; Value = '  -2.04%  ' },
    @{ Ref = 'E15'; Value = '  -0.61%  ' },
    @{ Ref = 'E16'; Value = '  +1.20%  ' },
    @{ Ref = 'E17'; Value = '  -1.30%  ' },
    @{ Ref = 'E18'; Value = '  +2.21%  ' },
    @{ Ref = 'E20'; Value = '  +1.51%  ' },
    @{ Ref = 'E21'; Value = '  -1.91%  ' },
    @{ Ref = 'E22'; Value = '  +0.09%  ' },
    @{ Ref = 'E23'; Value = '  -1.07%  ' },
    @{ Ref = 'E24'; Value = '  +5.01%  ' },
    @{ Ref = 'E25'; Value = '  -1.69%  ' },
    @{ Ref = 'E26'; Value = '  +0.36%  ' },
    @{ Ref = 'E27'; Value = '  -2.86%  ' },
    @{ Ref = 'E28'; Value = '  -5.62%  ' },
    @{ Ref = 'E29'; Value = '  -3.84%  ' },
    @{ Ref = 'E30'; Value = '  +0.05%  ' },
    @{ Ref = 'E31'; Value = '  -2.72%  ' },
    @{ Ref = 'E32'; Value = '  -2.43%  ' },
    @{ Ref = 'E33'; Value = '  -0.73%  ' },
    @{ Ref = 'E34'; Value = '  -3.40%  ' },
    @{ Ref = 'E35'; Value = '  +2.35%  ' },
    @{ Ref = 'B36'; Value = 'EthereumClassic' },
    @{ Ref = 'C36'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' },
    @{ Ref = 'E36'; Value = '  -0.08%  ' },
    @{ Ref = 'B37'; Value = 'PolygonEcosystemToken' },
    @{ Ref = 'C37'; Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol' },
    @{ Ref = 'E37'; Value = '  -1.97%  ' },
    @{ Ref = 'E38'; Value = '  +0.06%  ' },
    @{ Ref = 'E39'; Value = '  -2.98%  ' },
    @{ Ref = 'E40'; Value = '  -2.01%  ' },
    @{ Ref = 'E41'; Value = '  +0.05%  ' },
    @{ Ref = 'E42'; Value = '  -2.38%  ' },
    @{ Ref = 'E43'; Value = '  -0.87%  ' },
    @{ Ref = 'E44'; Value = '  +1.61%  ' },
    @{ Ref = 'E45'; Value = '  -0.39%  ' },
    @{ Ref = 'E46'; Value = '  -2.56%  ' },
    @{ Ref = 'E47'; Value = '  -2.33%  ' },
    @{ Ref = 'E48'; Value = '  +16.28%  ' },
    @{ Ref = 'E49'; Value = '  -2.87%  ' },
    @{ Ref = 'E50'; Value = '  +1.53%  ' },
    @{ Ref = 'E51'; Value = '  -4.15%  ' }
)

foreach ($item in $textCells) {
    $ws.Range($item.Ref).Value = $item.Value
}
